$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TC1 "Step" cell (D2): remove the birthday step and renumber the remaining steps.
$newStep = "1. Open browser`n2. Visit the website: ""https://demoqa.com/automation-practice-form""`n3. Input valid first name: 'Kha'`n4. Input valid last name: 'Hoai'`n5. Select Gender option: 'Male'`n6.Input valid mobile number: '9825467895'`n7. Input valid subject: 'Math' 'Chemistry'`n8. Click on Submit button"
$ws.Range("D2").Value = $newStep

# Row 2 auto-fits shorter now that it has one less wrapped line of text.
$ws.Rows.Item(2).RowHeight = 141.75

# The first picture is anchored absolutely, so its on-sheet position is fixed,
# but its cached "to" cell/row/offset needs to be refreshed against the new
# row 2 height. Touching Height forces that anchor cache to be recomputed.
$pic = $ws.Shapes.Item(1)
$pic.Height = $pic.Height()

# Update the saved view state: selection moves to D22 (scrolled further down
# the sheet), while row 1 stays frozen.
$ws.Range("D22").Select()
